# Weekly refresh of the "Bruselas (repollito) - Vega Modelo de Temuco" price log.
# The sheet currently holds 27 data rows (rows 2-28). This edit:
#   1) inserts 7 new rows (29-35) so the table grows to 34 data rows (rows 2-35),
#   2) fills in the constant columns for the new rows (same market/region/etc. as
#      every other row in this subset),
#   3) rewrites the per-row figures (Fecha, Volumen, Precio minimo/maximo/promedio,
#      Precio $/Kg) for rows 8-35 to their new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Make room for the 7 new data rows at the bottom of the table ---------
$ws.Range("A29:R35").EntireRow.Insert()

# --- 2) Constant columns for the newly inserted rows -------------------------
for ($r = 29; $r -le 35; $r++) {
    $ws.Cells.Item($r, 1).Value  = 10
    $ws.Cells.Item($r, 2).Value  = "Vega Modelo de Temuco"
    $ws.Cells.Item($r, 3).Value  = "La Araucanía"
    $ws.Cells.Item($r, 5).Value  = 9
    $ws.Cells.Item($r, 6).Value  = 100112035
    $ws.Cells.Item($r, 7).Value  = "Bruselas (repollito)"
    $ws.Cells.Item($r, 8).Value  = "Sin especificar"
    $ws.Cells.Item($r, 9).Value  = "Primera"
    $ws.Cells.Item($r, 14).Value = "`$/malla 10 kilos"
    $ws.Cells.Item($r, 15).Value = "Provincia de Quillota"
    $ws.Cells.Item($r, 17).Value = 10
    $ws.Cells.Item($r, 18).Value = "Hortaliza"
}

# --- 3) Per-row figures (rows 8-35) -------------------------------------------
$rows = @(
  @{ R = 8;  D = 44425; J = 30; K = 25000; L = 25000; M = 25000; P = 2500 },
  @{ R = 9;  D = 44392; J = 25; K = 24000; L = 24000; M = 24000; P = 2400 },
  @{ R = 10; D = 44411; J = 40; K = 25000; L = 25000; M = 25000; P = 2500 },
  @{ R = 11; D = 44379; J = 35; K = 22000; L = 22000; M = 22000; P = 2200 },
  @{ R = 12; D = 44356; J = 15; K = 24000; L = 24000; M = 24000; P = 2400 },
  @{ R = 13; D = 44396; J = 20; K = 25000; L = 25000; M = 25000; P = 2500 },
  @{ R = 14; D = 44384; J = 40; K = 25000; L = 25000; M = 25000; P = 2500 },
  @{ R = 15; D = 44349; J = 45; K = 24000; L = 24000; M = 24000; P = 2400 },
  @{ R = 16; D = 44385; J = 80; K = 25000; L = 25000; M = 25000; P = 2500 },
  @{ R = 17; D = 44427; J = 40; K = 25000; L = 25000; M = 25000; P = 2500 },
  @{ R = 18; D = 44413; J = 40; K = 25000; L = 25000; M = 25000; P = 2500 },
  @{ R = 19; D = 44421; J = 55; K = 25000; L = 25000; M = 25000; P = 2500 },
  @{ R = 20; D = 44400; J = 12; K = 24000; L = 24000; M = 24000; P = 2400 },
  @{ R = 21; D = 44426; J = 30; K = 25000; L = 25000; M = 25000; P = 2500 },
  @{ R = 22; D = 44390; J = 15; K = 25000; L = 25000; M = 25000; P = 2500 },
  @{ R = 23; D = 44354; J = 30; K = 24000; L = 24000; M = 24000; P = 2400 },
  @{ R = 24; D = 44410; J = 50; K = 25000; L = 25000; M = 25000; P = 2500 },
  @{ R = 25; D = 44412; J = 50; K = 25000; L = 25000; M = 25000; P = 2500 },
  @{ R = 26; D = 44371; J = 50; K = 25000; L = 25000; M = 25000; P = 2500 },
  @{ R = 27; D = 44405; J = 40; K = 25000; L = 25000; M = 25000; P = 2500 },
  @{ R = 28; D = 44350; J = 40; K = 24000; L = 25000; M = 24375; P = 2438 },
  @{ R = 29; D = 44389; J = 65; K = 25000; L = 25000; M = 25000; P = 2500 },
  @{ R = 30; D = 44417; J = 15; K = 25000; L = 25000; M = 25000; P = 2500 },
  @{ R = 31; D = 44419; J = 25; K = 25000; L = 25000; M = 25000; P = 2500 },
  @{ R = 32; D = 44420; J = 55; K = 25000; L = 25000; M = 25000; P = 2500 },
  @{ R = 33; D = 44382; J = 50; K = 25000; L = 25000; M = 25000; P = 2500 },
  @{ R = 34; D = 44355; J = 25; K = 23000; L = 24000; M = 23400; P = 2340 },
  @{ R = 35; D = 44376; J = 45; K = 23000; L = 23000; M = 23000; P = 2300 }
)

foreach ($row in $rows) {
    $r = $row.R
    $ws.Cells.Item($r, 4).Value  = $row.D   # D: Fecha
    $ws.Cells.Item($r, 10).Value = $row.J   # J: Volumen
    $ws.Cells.Item($r, 11).Value = $row.K   # K: Precio minimo
    $ws.Cells.Item($r, 12).Value = $row.L   # L: Precio maximo
    $ws.Cells.Item($r, 13).Value = $row.M   # M: Precio promedio ponderado
    $ws.Cells.Item($r, 16).Value = $row.P   # P: Precio $/Kg
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}

Write-Output "ok"
